# Insert one new weekly price record for "Ají" (Inferno / Primera) right
# before the existing row 489, shifting the remaining records down by one
# row (old row 489 -> new row 490, ..., old row 514 -> new row 515).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(489).Insert()

$ws.Cells.Item(489, 1).Value = 8
$ws.Cells.Item(489, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(489, 3).Value = "Coquimbo"
$ws.Cells.Item(489, 4).Value = 45147
$ws.Cells.Item(489, 5).Value = 4
$ws.Cells.Item(489, 6).Value = 100112021
$ws.Cells.Item(489, 7).Value = "Ají"
$ws.Cells.Item(489, 8).Value = "Inferno"
$ws.Cells.Item(489, 9).Value = "Primera"
$ws.Cells.Item(489, 10).Value = 400
$ws.Cells.Item(489, 11).Value = 15000
$ws.Cells.Item(489, 12).Value = 16000
$ws.Cells.Item(489, 13).Value = 15500
$ws.Cells.Item(489, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(489, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(489, 16).Value = 1550
$ws.Cells.Item(489, 17).Value = 10
$ws.Cells.Item(489, 18).Value = "Hortaliza"
